$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.773.20"
$ws.Range("E2").Value = "  +2.51%  "
$ws.Range("D3").Value = "1.876.31"
$ws.Range("E3").Value = "  +2.38%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.28%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "325.43"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").Value = "  +0.27%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4584"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.71%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3868"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07860"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.17%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.9942"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +3.46%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "1.881.28"
$ws.Range("E12").Value = "  +0.92%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "6.996"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.65%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.718"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.81%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.06953"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.21%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "88.51"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.37%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.00001006"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.26%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "16.87"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "28.780.91"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("E22").Value = "  -0.25%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "11.05"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.48%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.130"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").Value = "2.116.27"
$ws.Range("E25").Value = "  +1.82%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "153.40"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "19.23"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.41%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "5.805"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.53%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.972"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.31%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "119.15"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.41%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.09321"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.90%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.9192"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.26%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.305"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.86%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.342"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.57%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "3.325"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.05766"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.22%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.154"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +1.70%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.02075"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.86%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "7.711"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.51%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.5643"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +1.02%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.1790"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +1.74%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "9.934"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.34%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.07221"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -2.47%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "11.85"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +1.85%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.5304"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.67%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.152"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +2.19%  "
$ws.Range("E47").Value = "  -0.67%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "113.67"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("E49").Value = "  -0.21%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.410"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +3.89%  "
$ws.Range("E51").Value = "  +0.31%  "
